# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 01:06"

# --- Swap country names that changed rank order (mirrors shared-string reorder in the diff) ---
# Armenia (row 52) <-> Nigeria (row 53)
$ws.Range("A52").Value = "Nigeria"
$ws.Range("A53").Value = "Armenia"

# Cabo Verde(133), Republica del Chad(134), Principado de Andorra(135), Uruguay(136) -> rotate 134-136
$ws.Range("A134").Value = "Uruguay"
$ws.Range("A135").Value = "Republica del Chad"
$ws.Range("A136").Value = "Principado de Andorra"

# Fiyi (row 202) <-> Dominica (row 203)
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

# Santa Sede (row 208) <-> Islas Turcas y Caicos (row 209)
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("A209").Value = "Santa Sede"

# Papua Nueva Guinea (row 213) <-> Islas Virgenes Britanicas (row 214)
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("A214").Value = "Papua Nueva Guinea"

# --- Update numeric statistics ---
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2329027
$ws.Range("C4").Value = 31837
$ws.Range("D4").Value = 966595
$ws.Range("E4").Value = 1240479
$ws.Range("G4").Value = 546
$ws.Range("H4").Value = 121953

# Row 14
$ws.Range("B14").Value = 191216
$ws.Range("C14").Value = 556
$ws.Range("E14").Value = 7555

# Row 25
$ws.Range("B25").Value = 65633
$ws.Range("C25").Value = 2357
$ws.Range("D25").Value = 25499
$ws.Range("E25").Value = 38008
$ws.Range("G25").Value = 81
$ws.Range("H25").Value = 2126

# Row 52 (now Nigeria)
$ws.Range("B52").Value = 19808
$ws.Range("C52").Value = 661
$ws.Range("D52").Value = 6718
$ws.Range("E52").Value = 12584
$ws.Range("G52").Value = 19
$ws.Range("H52").Value = 506

# Row 53 (now Armenia)
$ws.Range("B53").Value = 19708
$ws.Range("C53").Value = 551
$ws.Range("D53").Value = 8854
$ws.Range("E53").Value = 10522
$ws.Range("G53").Value = 13
$ws.Range("H53").Value = 332

# Row 54
$ws.Range("B54").Value = 17799
$ws.Range("C54").Value = 59
$ws.Range("D54").Value = 16077
$ws.Range("E54").Value = 770
$ws.Range("G54").Value = 17
$ws.Range("H54").Value = 952

# Row 67
$ws.Range("B67").Value = 10448
$ws.Range("C67").Value = 42
$ws.Range("E67").Value = 2635

# Row 69
$ws.Range("B69").Value = 8742
$ws.Range("C69").Value = 16
$ws.Range("E69").Value = 360

# Row 134 (now Uruguay)
$ws.Range("B134").Value = 859
$ws.Range("C134").Value = 6
$ws.Range("D134").Value = 815
$ws.Range("E134").Value = 19
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 25

# Row 135 (now Republica del Chad)
$ws.Range("B135").Value = 858
$ws.Range("D135").Value = 746
$ws.Range("E135").Value = 38
$ws.Range("H135").Value = 74

# Row 136 (now Principado de Andorra)
$ws.Range("B136").Value = 855
$ws.Range("D136").Value = 792
$ws.Range("E136").Value = 11
$ws.Range("H136").Value = 52

# Row 157
$ws.Range("D157").Value = 327
$ws.Range("E157").Value = 22

# Row 208 (now Islas Turcas y Caicos)
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

# Row 209 (now Santa Sede)
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0

# Row 213 (now Islas Virgenes Britanicas)
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

# Row 214 (now Papua Nueva Guinea)
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
